$wb = $excel.ActiveWorkbook
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$weather = $wb.Worksheets.Add($null, $after)
$weather.Name = "weather"
$weather.Range("A1").Value = "Skycondition"
$weather.Range("B1").Value = "Partly Cloudy"
try {
    $weather.Columns.AutoFit()
    Write-Host "autofit ok"
} catch {
    Write-Host ("AUTOFIT ERROR: " + $_.Exception.Message)
}
